$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: G21 becomes a hard-coded value (no longer E21*D21 formula)
$ws.Range("G21").Value = 22.9

# Row 22: unit price changes
$ws.Range("E22").Value = 99.9

# Row 23: unit price changes (swapped back with row 22's old value)
$ws.Range("E23").Value = 17.9

# Row 24: unit price changes
$ws.Range("E24").Value = 16.9

# Row 25: unit price changes, and the cell's style reverts to the default
# (Calibri 11, black) matching the rest of the column instead of the
# Arial 10 "explanatory text" style it previously had.
$ws.Range("E25").Value = 17.9
$ws.Range("E25").Font.Name = "Calibri"
$ws.Range("E25").Font.Size = 11
$ws.Range("E25").Font.Color = 0

# Update the active selection to match the author's final cursor position
$ws.Range("G29").Select()
